# update application get api
#
# Appends three new enquiry rows (rows 6-8) to the worksheet, matching the
# data added upstream by the API update.
#
# Note: the source diff also widens the sheet's
#   <ignoredErrors><ignoredError numberStoredAsText="1" sqref="A1:L5"/></ignoredErrors>
# entry to sqref="A1:L8". That element is not reachable through the Excel
# COM object model (there is no Range/Worksheet/Application member that
# edits the persisted "Number Stored As Text" ignore-list), so it cannot be
# updated from this script. Everything else the diff describes (new rows
# 6-8 with their values, the growth of the sheet dimension to A1:L8, and
# the date styling on the new J/K cells) is applied below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 ---
$ws.Range("A6").Value = "3a98ce58-eda0-4aa8-990c-be9e0a3131e1"
$ws.Range("B6").Value = "Brooke Daniels"
$ws.Range("C6").Value = "Ochoa Merritt Inc"
$ws.Range("D6").Value = "Nisi dolor sapiente "
$ws.Range("E6").Value = "fyfu@mailinator.com"
$ws.Range("F6").Value = "university"
$ws.Range("G6").Value = "for-organization"
$ws.Range("H6").Value = "new"
$ws.Range("J6").Value = 45434.77850829861
$ws.Range("K6").Value = 45434.77850829861

# --- Row 7 ---
$ws.Range("A7").Value = "2139c963-0efb-4ecf-8f66-28129184f8ed"
$ws.Range("B7").Value = "Umair Rinde"
$ws.Range("C7").Value = "Clarke and Castro Trading"
$ws.Range("D7").Value = "Officiis voluptates "
$ws.Range("E7").Value = "rindeumair@gmail.com"
$ws.Range("F7").Value = "online"
$ws.Range("G7").Value = "for-organization"
$ws.Range("H7").Value = "new"
$ws.Range("J7").Value = 45434.77882318287
$ws.Range("K7").Value = 45434.77882318287

# --- Row 8 ---
$ws.Range("A8").Value = "c226674c-7f3c-436a-ae2b-f3f7391b952b"
$ws.Range("B8").Value = "Umair Rinde"
$ws.Range("C8").Value = "Goff Walsh Plc"
$ws.Range("D8").Value = "Officia eligendi dic"
$ws.Range("E8").Value = "rindeumair@gmail.com"
$ws.Range("F8").Value = "online"
$ws.Range("G8").Value = "for-partnership"
$ws.Range("H8").Value = "new"
$ws.Range("J8").Value = 45434.77936392361
$ws.Range("K8").Value = 45434.77936392361

# Copy the date/time number format (the style index used by the existing
# createdAt/updatedAt columns) onto the new J/K cells so the new rows are
# styled exactly like the existing ones (this reuses the same style index
# instead of creating new, duplicate styles).
$ws.Range("J2:K2").Copy() | Out-Null
$ws.Range("J6:K6").PasteSpecial(-4122) | Out-Null
$ws.Range("J7:K7").PasteSpecial(-4122) | Out-Null
$ws.Range("J8:K8").PasteSpecial(-4122) | Out-Null

# PasteSpecial only touches formatting, so make sure the values are still
# the intended ones.
$ws.Range("J6").Value = 45434.77850829861
$ws.Range("K6").Value = 45434.77850829861
$ws.Range("J7").Value = 45434.77882318287
$ws.Range("K7").Value = 45434.77882318287
$ws.Range("J8").Value = 45434.77936392361
$ws.Range("K8").Value = 45434.77936392361
